$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row in column A (Beteckning) to know the data extent.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C ("Förändrad") holds a date serial number for every data row (2..lastRow).
# Update each cell's value from 45205 to 45206, preserving existing formatting.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
